$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: age value changes from 43 to 45
$ws.Range("A5").Value = 45

# Previously-empty row 6 gets the record that used to be age 43 (male, zip 81931)
$ws.Range("A6").Value = 43
$ws.Range("B6").Value = "male"
$ws.Range("C6").Value = 81931

# Previously-empty row 16 gets a new record: age 70, male, zip 81931
$ws.Range("A16").Value = 70
$ws.Range("B16").Value = "male"
$ws.Range("C16").Value = 81931

# The record that used to live in row 17 (45, male, 81931) now lives in row 18
$ws.Range("A18").Value = 45
$ws.Range("B18").Value = "male"
$ws.Range("C18").Value = 81931

# Row 17 is now empty again
$ws.Range("A17:C17").ClearContents()

# Update the active selection to reflect where the user ended up working
$ws.Range("C6").Select()
